$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("issues")

# The "Manual Testing Multi dataset fitting" row (row 10) was removed entirely;
# deleting the row shifts every row below it up by one and Excel automatically
# renumbers/garbage-collects the shared-string table and recalculates the
# dependent COUNTIF formulas on the "assignees" sheet.
$ws.Rows(10).Delete()

# The "Manual Testing Mantid Basics Course" row (now row 17 after the delete)
# gained an extra bullet point about building docs-qtassistant/docs-qthelp.
$ws.Range("B17").Value = "* Build a docs-html target of Mantid`n* Build docs-qtassistant and docs-qthelp so that autogenerated plots are available`n* Check against the online documentation (http://docs.mantidproject.org/nightly/tutorials/mantid_basic_course/index.html#mantid-basic-course)`n* Open up the basic course (docs/html/tutorials/mantid_basic_course/index.html)`n* Check that the pages in there make sense"

# Restore the selection used by the author after editing.
$ws.Range("C17").Select()
